$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(20, 8).Value = 2149.4285
$ws.Cells.Item(20, 9).Value = 402.5
$ws.Cells.Item(20, 10).Value = 12631
$ws.Cells.Item(20, 11).Value = 402.5
$ws.Cells.Item(20, 12).Value = 12631
$ws.Cells.Item(20, 13).Value = -172.5
$ws.Cells.Item(20, 14).Value = -13091

$ws.Cells.Item(35, 8).Value = 2149.4285
$ws.Cells.Item(35, 9).Value = 402.5
$ws.Cells.Item(35, 10).Value = 12631
$ws.Cells.Item(35, 11).Value = 402.5
$ws.Cells.Item(35, 12).Value = 12631
$ws.Cells.Item(35, 13).Value = -23.5
$ws.Cells.Item(35, 14).Value = -13389

$ws.Cells.Item(40, 8).Value = 4287.04
$ws.Cells.Item(40, 9).Value = 3688.0715
$ws.Cells.Item(40, 11).Value = 3688.0715
$ws.Cells.Item(40, 13).Value = -3513.0715

$ws.Cells.Item(51, 8).Value = 0
$ws.Cells.Item(51, 10).Value = 0
$ws.Cells.Item(51, 12).Value = 0
$ws.Cells.Item(51, 14).ClearContents()

$ws.Cells.Item(57, 8).Value = 48890
$ws.Cells.Item(57, 9).Value = 47000
$ws.Cells.Item(57, 10).Value = 50780
$ws.Cells.Item(57, 11).Value = 141000
$ws.Cells.Item(57, 12).Value = 152340
$ws.Cells.Item(57, 13).Value = -140501
$ws.Cells.Item(57, 14).Value = -153338

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2526.8635
$ws.Cells.Item(2, 9).Value = 2079.1875
$ws.Cells.Item(2, 10).Value = 3720.6667
$ws.Cells.Item(2, 11).Value = 2079.1875
$ws.Cells.Item(2, 12).Value = 3720.6667
$ws.Cells.Item(2, 13).Value = -1966.1875
$ws.Cells.Item(2, 14).Value = -3946.6667

$ws.Cells.Item(32, 8).Value = 5550.56
$ws.Cells.Item(32, 9).Value = 4076.6956
$ws.Cells.Item(32, 11).Value = 4076.6956
$ws.Cells.Item(32, 13).Value = -3789.6956

$ws.Cells.Item(116, 8).Value = 2526.8635
$ws.Cells.Item(116, 9).Value = 2079.1875
$ws.Cells.Item(116, 10).Value = 3720.6667
$ws.Cells.Item(116, 11).Value = 2079.1875
$ws.Cells.Item(116, 12).Value = 3720.6667
$ws.Cells.Item(116, 13).Value = 214.8125
$ws.Cells.Item(116, 14).Value = -8308.6667

$ws.Cells.Item(122, 8).Value = 2746.0833
$ws.Cells.Item(122, 9).Value = 1751.125
$ws.Cells.Item(122, 10).Value = 4736
$ws.Cells.Item(122, 11).Value = 5253.375
$ws.Cells.Item(122, 12).Value = 14208
$ws.Cells.Item(122, 13).Value = -2803.375
$ws.Cells.Item(122, 14).Value = -19108

$ws.Cells.Item(132, 8).Value = 897.5
$ws.Cells.Item(132, 9).Value = 897.5
$ws.Cells.Item(132, 11).Value = 2692.5
$ws.Cells.Item(132, 13).Value = -162.5

$ws.Cells.Item(133, 8).Value = 55000
$ws.Cells.Item(133, 10).Value = 55000
$ws.Cells.Item(133, 12).Value = 55000
$ws.Cells.Item(133, 14).Value = -60060

$ws.Cells.Item(135, 8).Value = 100428.5
$ws.Cells.Item(135, 10).Value = 100428.5
$ws.Cells.Item(135, 14).Value = -110568.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2526.8635
$ws.Cells.Item(3, 9).Value = 2079.1875
$ws.Cells.Item(3, 10).Value = 3720.6667
$ws.Cells.Item(3, 11).Value = 2079.1875
$ws.Cells.Item(3, 12).Value = 3720.6667
$ws.Cells.Item(3, 13).Value = -1965.1875
$ws.Cells.Item(3, 14).Value = -3948.6667

$ws.Cells.Item(86, 8).Value = 2582.3572
$ws.Cells.Item(86, 9).Value = 1294.95
$ws.Cells.Item(86, 11).Value = 1294.95
$ws.Cells.Item(86, 13).Value = -171.95

$ws.Cells.Item(89, 8).Value = 2582.3572
$ws.Cells.Item(89, 9).Value = 1294.95
$ws.Cells.Item(89, 11).Value = 6474.75
$ws.Cells.Item(89, 13).Value = -858.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1376.125
$ws.Cells.Item(16, 9).Value = 1252.25
$ws.Cells.Item(16, 10).Value = 1500
$ws.Cells.Item(16, 11).Value = 1252.25
$ws.Cells.Item(16, 12).Value = 1500
$ws.Cells.Item(16, 13).Value = -965.25
$ws.Cells.Item(16, 14).Value = -2074

$ws.Cells.Item(22, 8).Value = 2731
$ws.Cells.Item(22, 9).Value = 1619.6
$ws.Cells.Item(22, 11).Value = 1619.6
$ws.Cells.Item(22, 13).Value = -1269.6

$ws.Cells.Item(99, 8).Value = 3128.5
$ws.Cells.Item(99, 9).Value = 2700
$ws.Cells.Item(99, 11).Value = 2700
$ws.Cells.Item(99, 13).Value = -1202

$ws.Cells.Item(113, 8).Value = 1376.125
$ws.Cells.Item(113, 9).Value = 1252.25
$ws.Cells.Item(113, 10).Value = 1500
$ws.Cells.Item(113, 11).Value = 1252.25
$ws.Cells.Item(113, 12).Value = 1500
$ws.Cells.Item(113, 13).Value = 917.75
$ws.Cells.Item(113, 14).Value = -5840

$ws.Cells.Item(122, 8).Value = 1062.2858
$ws.Cells.Item(122, 9).Value = 1006.5
$ws.Cells.Item(122, 10).Value = 1201.75
$ws.Cells.Item(122, 11).Value = 3019.5
$ws.Cells.Item(122, 12).Value = 3605.25
$ws.Cells.Item(122, 13).Value = -569.5
$ws.Cells.Item(122, 14).Value = -8505.25

$ws.Cells.Item(126, 8).Value = 3128.5
$ws.Cells.Item(126, 9).Value = 2700
$ws.Cells.Item(126, 11).Value = 8100
$ws.Cells.Item(126, 13).Value = -5630

$ws.Cells.Item(132, 8).Value = 4813.769
$ws.Cells.Item(132, 9).Value = 4179.8335
$ws.Cells.Item(132, 10).Value = 5357.143
$ws.Cells.Item(132, 11).Value = 12539.5005
$ws.Cells.Item(132, 12).Value = 16071.429
$ws.Cells.Item(132, 13).Value = -10009.5005
$ws.Cells.Item(132, 14).Value = -21131.429

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 19.166666
$ws.Cells.Item(2, 9).Value = 19.2
$ws.Cells.Item(2, 10).Value = 19.1
$ws.Cells.Item(2, 11).Value = 115.2
$ws.Cells.Item(2, 12).Value = 114.6
$ws.Cells.Item(2, 13).Value = -2.199999999999989
$ws.Cells.Item(2, 14).Value = -340.6

$ws.Cells.Item(7, 8).Value = 44.75
$ws.Cells.Item(7, 9).Value = 25
$ws.Cells.Item(7, 10).Value = 56.6
$ws.Cells.Item(7, 11).Value = 75
$ws.Cells.Item(7, 12).Value = 169.8
$ws.Cells.Item(7, 13).Value = 37
$ws.Cells.Item(7, 14).Value = -393.8

$ws.Cells.Item(34, 8).Value = 1778.4445
$ws.Cells.Item(34, 10).Value = 3498.2856
$ws.Cells.Item(34, 12).Value = 10494.8568
$ws.Cells.Item(34, 14).Value = -10662.8568

$ws.Cells.Item(92, 8).Value = 50001.5
$ws.Cells.Item(92, 10).Value = 90003
$ws.Cells.Item(92, 12).Value = 270009
$ws.Cells.Item(92, 14).Value = -272505

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 54.5
$ws.Cells.Item(2, 9).Value = 47.5
$ws.Cells.Item(2, 11).Value = 47.5
$ws.Cells.Item(2, 13).Value = 65.5

$ws.Cells.Item(6, 8).Value = 83.333336
$ws.Cells.Item(6, 10).Value = 83.333336
$ws.Cells.Item(6, 12).Value = 83.333336
$ws.Cells.Item(6, 14).Value = -309.333336

$ws.Cells.Item(16, 8).Value = 83.333336
$ws.Cells.Item(16, 10).Value = 83.333336
$ws.Cells.Item(16, 12).Value = 83.333336
$ws.Cells.Item(16, 14).Value = -583.333336

$ws.Cells.Item(59, 8).Value = 0
$ws.Cells.Item(59, 9).Value = 0
$ws.Cells.Item(59, 10).Value = 0
$ws.Cells.Item(59, 11).Value = 0
$ws.Cells.Item(59, 12).Value = 0
$ws.Cells.Item(59, 13).ClearContents()
$ws.Cells.Item(59, 14).ClearContents()

$ws.Cells.Item(80, 8).Value = 4933.3335
$ws.Cells.Item(80, 9).Value = 4933.3335
$ws.Cells.Item(80, 11).Value = 4933.3335
$ws.Cells.Item(80, 13).Value = -3935.3335

$ws.Cells.Item(83, 8).Value = 4933.3335
$ws.Cells.Item(83, 9).Value = 4933.3335
$ws.Cells.Item(83, 11).Value = 24666.6675
$ws.Cells.Item(83, 13).Value = -19674.6675

$ws.Cells.Item(122, 8).Value = 2331.5293
$ws.Cells.Item(122, 9).Value = 1977.5625
$ws.Cells.Item(122, 10).Value = 7995
$ws.Cells.Item(122, 11).Value = 5932.6875
$ws.Cells.Item(122, 12).Value = 23985
$ws.Cells.Item(122, 13).Value = -3482.6875
$ws.Cells.Item(122, 14).Value = -28885

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(35, 8).Value = 2562.7778
$ws.Cells.Item(35, 9).Value = 1498.75
$ws.Cells.Item(35, 10).Value = 3414
$ws.Cells.Item(35, 11).Value = 1498.75
$ws.Cells.Item(35, 12).Value = 3414
$ws.Cells.Item(35, 13).Value = -1162.75
$ws.Cells.Item(35, 14).Value = -4086

$ws.Cells.Item(61, 8).Value = 3772.353
$ws.Cells.Item(61, 10).Value = 5597.75
$ws.Cells.Item(61, 12).Value = 5597.75
$ws.Cells.Item(61, 14).Value = -6001.75

$ws.Cells.Item(93, 8).Value = 2234.7144
$ws.Cells.Item(93, 9).Value = 2273.8333
$ws.Cells.Item(93, 10).Value = 2000
$ws.Cells.Item(93, 11).Value = 2273.8333
$ws.Cells.Item(93, 12).Value = 2000
$ws.Cells.Item(93, 13).Value = -1025.8333
$ws.Cells.Item(93, 14).Value = -4496

$ws.Cells.Item(113, 8).Value = 3772.353
$ws.Cells.Item(113, 10).Value = 5597.75
$ws.Cells.Item(113, 12).Value = 5597.75
$ws.Cells.Item(113, 14).Value = -9937.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(5, 8).Value = 9700500
$ws.Cells.Item(5, 10).Value = 5400000
$ws.Cells.Item(5, 12).Value = 5400000
$ws.Cells.Item(5, 14).Value = -5400224

$ws.Cells.Item(113, 8).Value = 1017.125
$ws.Cells.Item(113, 9).Value = 1248
$ws.Cells.Item(113, 11).Value = 3744
$ws.Cells.Item(113, 13).Value = -1574

$ws.Cells.Item(126, 8).Value = 4487.8335
$ws.Cells.Item(126, 9).Value = 890.8
$ws.Cells.Item(126, 11).Value = 2672.4
$ws.Cells.Item(126, 13).Value = -202.3999999999996

$ws.Cells.Item(136, 8).Value = 2915
$ws.Cells.Item(136, 9).Value = 1221.6666
$ws.Cells.Item(136, 10).Value = 7995
$ws.Cells.Item(136, 11).Value = 3664.9998
$ws.Cells.Item(136, 12).Value = 23985
$ws.Cells.Item(136, 13).Value = -1114.9998
$ws.Cells.Item(136, 14).Value = -29085
